$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D values: "MSI" header-ish label (row10) + MSI numbers for the
# data rows below it. D10 keeps the existing (hyperlink-font) style already
# applied to that cell; D11-D14 get the plain text style used by the rest
# of the row (same as column A/B/C in those rows).

$ws.Range("D10").Value = "MSI"

$ws.Range("D11").Value = "732111198172294"
$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"

$ws.Range("D12").Value = "732111198172293"
$ws.Range("D12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"

$ws.Range("D13").Value = "732111198172294"
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"

$ws.Range("D14").Value = "732111198172293"
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"

# Move the active selection, matching the saved cursor position.
$ws.Range("D19").Select() | Out-Null
